$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-12 17:18:26"
$ws.Range("K2").Value = "7.1 MJ/m2"
$ws.Range("E3").Value = "2026-02-12 17:18:29"
$ws.Range("H3").Value = "'71%"
$ws.Range("K3").Value = "8.5 MJ/m2"
$ws.Range("O3").Value = "-3.1 °C"
$ws.Range("E4").Value = "2026-02-12 17:18:31"
$ws.Range("J4").Value = "998.1 hPa"
$ws.Range("K4").Value = "13.4 MJ/m2"
$ws.Range("E5").Value = "2026-02-12 17:18:34"
$ws.Range("K5").Value = "6.9 MJ/m2"
$ws.Range("E6").Value = "2026-02-12 17:18:36"
$ws.Range("J6").Value = "998.0 hPa"
$ws.Range("K6").Value = "13.8 MJ/m2"
$ws.Range("O6").Value = "16.2 °C"
$ws.Range("E7").Value = "2026-02-12 17:18:38"
$ws.Range("H7").Value = "'34%"
$ws.Range("J7").Value = "1000.9 hPa"
$ws.Range("K7").Value = "14.0 MJ/m2"
$ws.Range("E8").Value = "2026-02-12 17:18:41"
$ws.Range("J8").Value = "1000.2 hPa"
$ws.Range("K8").Value = "14.4 MJ/m2"
$ws.Range("E9").Value = "2026-02-12 17:18:43"
$ws.Range("H9").Value = "'71%"
$ws.Range("K9").Value = "13.4 MJ/m2"
$ws.Range("E10").Value = "2026-02-12 17:18:46"
$ws.Range("K10").Value = "13.6 MJ/m2"
$ws.Range("O10").Value = "15.2 °C"
$ws.Range("E11").Value = "2026-02-12 17:18:48"
$ws.Range("E12").Value = "2026-02-12 17:18:51"
$ws.Range("H12").Value = "'77%"
$ws.Range("O12").Value = "12.5 °C"
$ws.Range("E13").Value = "2026-02-12 17:18:53"
$ws.Range("H13").Value = "'51%"
$ws.Range("J13").Value = "1000.8 hPa"
$ws.Range("K13").Value = "13.4 MJ/m2"
$ws.Range("E14").Value = "2026-02-12 17:18:56"
$ws.Range("K14").Value = "14.0 MJ/m2"
$ws.Range("E15").Value = "2026-02-12 17:18:58"
$ws.Range("E16").Value = "2026-02-12 17:19:01"
$ws.Range("K16").Value = "10.5 MJ/m2"
$ws.Range("E17").Value = "2026-02-12 17:19:03"
$ws.Range("K17").Value = "14.9 MJ/m2"
$ws.Range("E18").Value = "2026-02-12 17:19:05"
$ws.Range("J18").Value = "998.4 hPa"
$ws.Range("K18").Value = "13.8 MJ/m2"
$ws.Range("E19").Value = "2026-02-12 17:19:08"
$ws.Range("K19").Value = "13.7 MJ/m2"
$ws.Range("O19").Value = "8.4 °C"
$ws.Range("E20").Value = "2026-02-12 17:19:11"
$ws.Range("H20").Value = "'83%"
$ws.Range("K20").Value = "14.7 MJ/m2"
$ws.Range("E21").Value = "2026-02-12 17:19:13"
$ws.Range("H21").Value = "'49%"
$ws.Range("J21").Value = "1001.2 hPa"
$ws.Range("K21").Value = "13.2 MJ/m2"
$ws.Range("E22").Value = "2026-02-12 17:19:16"
$ws.Range("H22").Value = "'76%"
$ws.Range("K22").Value = "15.3 MJ/m2"
$ws.Range("O22").Value = "-5.5 °C"
$ws.Range("E23").Value = "2026-02-12 17:19:18"
$ws.Range("H23").Value = "'71%"
$ws.Range("K23").Value = "10.8 MJ/m2"
$ws.Range("E24").Value = "2026-02-12 17:19:21"
$ws.Range("J24").Value = "1006.0 hPa"
$ws.Range("K24").Value = "14.2 MJ/m2"
$ws.Range("E25").Value = "2026-02-12 17:19:23"
$ws.Range("H25").Value = "'61%"
$ws.Range("I25").Value = "2.0 mm"
$ws.Range("K25").Value = "14.1 MJ/m2"
$ws.Range("E26").Value = "2026-02-12 17:19:26"
$ws.Range("J26").Value = "997.2 hPa"
$ws.Range("K26").Value = "13.8 MJ/m2"
$ws.Range("O26").Value = "6.3 °C"
$ws.Range("E27").Value = "2026-02-12 17:19:28"
$ws.Range("H27").Value = "'64%"
$ws.Range("K27").Value = "13.8 MJ/m2"
$ws.Range("L27").Value = "59.0 km/h - 333º 16:31 TU"
$ws.Range("E28").Value = "2026-02-12 17:19:31"
$ws.Range("J28").Value = "997.6 hPa"
$ws.Range("K28").Value = "13.3 MJ/m2"
$ws.Range("E29").Value = "2026-02-12 17:19:33"
$ws.Range("K29").Value = "13.6 MJ/m2"
$ws.Range("O29").Value = "15.4 °C"
$ws.Range("E30").Value = "2026-02-12 17:19:36"
$ws.Range("H30").Value = "'68%"
$ws.Range("J30").Value = "998.2 hPa"
$ws.Range("K30").Value = "13.7 MJ/m2"
$ws.Range("O30").Value = "12.7 °C"
$ws.Range("E31").Value = "2026-02-12 17:19:38"
$ws.Range("J31").Value = "997.6 hPa"
$ws.Range("O31").Value = "14.6 °C"
$ws.Range("E32").Value = "2026-02-12 17:19:41"
$ws.Range("K32").Value = "13.9 MJ/m2"
$ws.Range("E33").Value = "2026-02-12 17:19:43"
$ws.Range("J33").Value = "1000.5 hPa"
$ws.Range("E34").Value = "2026-02-12 17:19:46"
$ws.Range("K34").Value = "14.5 MJ/m2"
$ws.Range("O34").Value = "0.5 °C"
$ws.Range("E35").Value = "2026-02-12 17:19:48"
$ws.Range("O35").Value = "8.0 °C"
$ws.Range("E36").Value = "2026-02-12 17:19:51"
$ws.Range("H36").Value = "'65%"
$ws.Range("J36").Value = "998.6 hPa"
$ws.Range("K36").Value = "13.7 MJ/m2"
$ws.Range("E37").Value = "2026-02-12 17:19:53"
$ws.Range("H37").Value = "'47%"
$ws.Range("J37").Value = "998.9 hPa"
$ws.Range("E38").Value = "2026-02-12 17:19:56"
$ws.Range("E39").Value = "2026-02-12 17:19:58"
$ws.Range("E40").Value = "2026-02-12 17:20:01"
$ws.Range("J40").Value = "1002.0 hPa"
$ws.Range("O40").Value = "10.2 °C"
$ws.Range("E41").Value = "2026-02-12 17:20:03"
$ws.Range("J41").Value = "1005.0 hPa"
$ws.Range("K41").Value = "14.2 MJ/m2"
$ws.Range("O41").Value = "17.4 °C"
$ws.Range("E42").Value = "2026-02-12 17:20:06"
$ws.Range("H42").Value = "'60%"
$ws.Range("E43").Value = "2026-02-12 17:20:08"
$ws.Range("K43").Value = "13.8 MJ/m2"
$ws.Range("E44").Value = "2026-02-12 17:20:10"
$ws.Range("E45").Value = "2026-02-12 17:20:13"
$ws.Range("J45").Value = "1004.0 hPa"
$ws.Range("K45").Value = "6.8 MJ/m2"
$ws.Range("E46").Value = "2026-02-12 17:20:15"
$ws.Range("J46").Value = "1006.8 hPa"
$ws.Range("K46").Value = "13.4 MJ/m2"
